# "many changes in framework"
# The Password column values are now stored base64-encoded instead of
# plain text, the Password column is a bit wider, and the last
# selection/active cell moved to B5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Base64-encode the existing plaintext passwords (Retrix123#) in column B.
$ws.Range("B2").Value = "UmV0cml4MTIzIw"
$ws.Range("B3").Value = "UmV0cml4MTIzIw"

# Widen the Password column a little to fit the longer encoded text.
$ws.Columns("B").ColumnWidth = 16.33

# Leave the cursor on B5, as last left by the editor.
$ws.Range("B5").Select()
